$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1253.7273
$ws.Range("J2").Value = 1311
$ws.Range("L2").Value = 1311
$ws.Range("N2").Value = -1537
$ws.Range("H17").Value = 2500
$ws.Range("J17").Value = 2500
$ws.Range("L17").Value = 7500
$ws.Range("N17").Value = -7836
$ws.Range("H19").Value = 9248.333000000001
$ws.Range("I19").Value = 9249
$ws.Range("J19").Value = 9248.200000000001
$ws.Range("K19").Value = 9249
$ws.Range("L19").Value = 9248.200000000001
$ws.Range("M19").Value = -9074
$ws.Range("N19").Value = -9598.200000000001
$ws.Range("H32").Value = 8835.3125
$ws.Range("J32").Value = 5225
$ws.Range("L32").Value = 5225
$ws.Range("N32").Value = -5877
$ws.Range("H40").Value = 4450
$ws.Range("I40").Value = 3675
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 3675
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -3500
$ws.Range("N40").Value = -6350
$ws.Range("H53").Value = 523.9167
$ws.Range("J53").Value = 479.875
$ws.Range("L53").Value = 479.875
$ws.Range("N53").Value = -1753.875
$ws.Range("H55").Value = 16.25
$ws.Range("I55").Value = 11.666667
$ws.Range("K55").Value = 11.666667
$ws.Range("M55").Value = 202.333333
$ws.Range("H62").Value = 7665.905
$ws.Range("I62").Value = 6905.3125
$ws.Range("K62").Value = 6905.3125
$ws.Range("M62").Value = -6281.3125
$ws.Range("H65").Value = 7665.905
$ws.Range("I65").Value = 6905.3125
$ws.Range("K65").Value = 34526.5625
$ws.Range("M65").Value = -31406.5625
$ws.Range("H76").Value = 4064.7
$ws.Range("J76").Value = 4333.3335
$ws.Range("L76").Value = 4333.3335
$ws.Range("N76").Value = -4963.3335
$ws.Range("H79").Value = 4064.7
$ws.Range("J79").Value = 4333.3335
$ws.Range("L79").Value = 4333.3335
$ws.Range("N79").Value = -6517.3335
$ws.Range("H88").Value = 18267.666
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 18267.666
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H96").Value = 7143608
$ws.Range("I96").Value = 8928860
$ws.Range("K96").Value = 26786580
$ws.Range("M96").Value = -26785207
$ws.Range("H107").Value = 1374.3462
$ws.Range("I107").Value = 1210.05
$ws.Range("K107").Value = 1210.05
$ws.Range("M107").Value = 709.95
$ws.Range("H113").Value = 4337.8823
$ws.Range("I113").Value = 3815.6667
$ws.Range("K113").Value = 3815.6667
$ws.Range("M113").Value = -561.6667000000002
$ws.Range("H116").Value = 26529.072
$ws.Range("I116").Value = 46726.2
$ws.Range("J116").Value = 15308.444
$ws.Range("K116").Value = 46726.2
$ws.Range("L116").Value = 15308.444
$ws.Range("M116").Value = -43284.2
$ws.Range("N116").Value = -22192.444
$ws.Range("H125").Value = 22668.857
$ws.Range("I125").Value = 75616
$ws.Range("J125").Value = 1490
$ws.Range("K125").Value = 680544
$ws.Range("L125").Value = 13410
$ws.Range("M125").Value = -678084
$ws.Range("N125").Value = -18330
$ws.Range("H129").Value = 1898.5454
$ws.Range("I129").Value = 1341.1428
$ws.Range("J129").Value = 2874
$ws.Range("K129").Value = 4023.4284
$ws.Range("L129").Value = 8622
$ws.Range("M129").Value = 976.5715999999998
$ws.Range("N129").Value = -18622
$ws.Range("H131").Value = 3087.5557
$ws.Range("I131").Value = 2738.4666
$ws.Range("K131").Value = 8215.399800000001
$ws.Range("M131").Value = -3175.399800000001
$ws.Range("H132").Value = 4103230
$ws.Range("I132").Value = 5522770
$ws.Range("K132").Value = 16568310
$ws.Range("M132").Value = -16565780
$ws.Range("H137").Value = 15030.074
$ws.Range("I137").Value = 19796.611
$ws.Range("J137").Value = 5497
$ws.Range("K137").Value = 59389.833
$ws.Range("L137").Value = 16491
$ws.Range("M137").Value = -56839.833
$ws.Range("N137").Value = -21591
$ws.Range("H138").Value = 20445.896
$ws.Range("I138").Value = 2169.6
$ws.Range("J138").Value = 40027.645
$ws.Range("K138").Value = 6508.799999999999
$ws.Range("L138").Value = 120082.935
$ws.Range("M138").Value = -1368.799999999999
$ws.Range("N138").Value = -130362.935
$ws.Range("H141").Value = 1249
$ws.Range("I141").Value = 999
$ws.Range("J141").Value = 1499
$ws.Range("K141").Value = 2997
$ws.Range("L141").Value = 4497
$ws.Range("M141").Value = 2183
$ws.Range("N141").Value = -14857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2571.7727
$ws.Range("I2").Value = 2801.1765
$ws.Range("K2").Value = 2801.1765
$ws.Range("M2").Value = -2688.1765
$ws.Range("H6").Value = 499050.75
$ws.Range("I6").Value = 400
$ws.Range("K6").Value = 400
$ws.Range("M6").Value = -227
$ws.Range("H32").Value = 15049.798
$ws.Range("I32").Value = 15579.16
$ws.Range("J32").Value = 5124.25
$ws.Range("K32").Value = 15579.16
$ws.Range("L32").Value = 5124.25
$ws.Range("M32").Value = -15292.16
$ws.Range("N32").Value = -5698.25
$ws.Range("H45").Value = 3454.8262
$ws.Range("I45").Value = 2208.0715
$ws.Range("J45").Value = 5394.222
$ws.Range("K45").Value = 2208.0715
$ws.Range("L45").Value = 5394.222
$ws.Range("M45").Value = -1831.0715
$ws.Range("N45").Value = -6148.222
$ws.Range("H61").Value = 6213.9546
$ws.Range("I61").Value = 1585.35
$ws.Range("K61").Value = 1585.35
$ws.Range("M61").Value = -1373.35
$ws.Range("H74").Value = 306080.16
$ws.Range("I74").Value = 501149.5
$ws.Range("J74").Value = 13476.125
$ws.Range("K74").Value = 501149.5
$ws.Range("L74").Value = 13476.125
$ws.Range("M74").Value = -500275.5
$ws.Range("N74").Value = -15224.125
$ws.Range("H77").Value = 306080.16
$ws.Range("I77").Value = 501149.5
$ws.Range("J77").Value = 13476.125
$ws.Range("K77").Value = 2505747.5
$ws.Range("L77").Value = 67380.625
$ws.Range("M77").Value = -2501379.5
$ws.Range("N77").Value = -76116.625
$ws.Range("H97").Value = 2309.389
$ws.Range("I97").Value = 1698
$ws.Range("K97").Value = 1698
$ws.Range("M97").Value = -1202
$ws.Range("H109").Value = 163333.33
$ws.Range("J109").Value = 163333.33
$ws.Range("L109").Value = 163333.33
$ws.Range("N109").Value = -166107.33
$ws.Range("H110").Value = 1436.5454
$ws.Range("I110").Value = 1054.2667
$ws.Range("K110").Value = 1054.2667
$ws.Range("M110").Value = 990.7333000000001
$ws.Range("H112").Value = 41519.57
$ws.Range("J112").Value = 41519.57
$ws.Range("L112").Value = 41519.57
$ws.Range("N112").Value = -44473.57
$ws.Range("H116").Value = 2571.7727
$ws.Range("I116").Value = 2801.1765
$ws.Range("K116").Value = 2801.1765
$ws.Range("M116").Value = -507.1765
$ws.Range("H122").Value = 1536.3914
$ws.Range("I122").Value = 1411.0555
$ws.Range("J122").Value = 1987.6
$ws.Range("K122").Value = 4233.166499999999
$ws.Range("L122").Value = 5962.799999999999
$ws.Range("M122").Value = -1783.166499999999
$ws.Range("N122").Value = -10862.8
$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030
$ws.Range("H136").Value = 6213.9546
$ws.Range("I136").Value = 1585.35
$ws.Range("K136").Value = 4756.049999999999
$ws.Range("M136").Value = -2206.049999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2571.7727
$ws.Range("I3").Value = 2801.1765
$ws.Range("K3").Value = 2801.1765
$ws.Range("M3").Value = -2687.1765
$ws.Range("H20").Value = 27686.45
$ws.Range("I20").Value = 45207.75
$ws.Range("J20").Value = 1404.5
$ws.Range("K20").Value = 45207.75
$ws.Range("L20").Value = 1404.5
$ws.Range("M20").Value = -44960.75
$ws.Range("N20").Value = -1898.5
$ws.Range("H22").Value = 840.6667
$ws.Range("I22").Value = 888.25
$ws.Range("J22").Value = 745.5
$ws.Range("K22").Value = 888.25
$ws.Range("L22").Value = 745.5
$ws.Range("M22").Value = -715.25
$ws.Range("N22").Value = -1091.5
$ws.Range("H74").Value = 61000
$ws.Range("J74").Value = 61000
$ws.Range("L74").Value = 61000
$ws.Range("N74").Value = -62872
$ws.Range("H77").Value = 61000
$ws.Range("J77").Value = 61000
$ws.Range("L77").Value = 183000
$ws.Range("N77").Value = -192360
$ws.Range("H82").Value = 23361.777
$ws.Range("H85").Value = 23361.777
$ws.Range("H86").Value = 1989.8
$ws.Range("I86").Value = 1989.8
$ws.Range("K86").Value = 1989.8
$ws.Range("M86").Value = -866.8
$ws.Range("H89").Value = 1989.8
$ws.Range("I89").Value = 1989.8
$ws.Range("K89").Value = 9949
$ws.Range("M89").Value = -4333
$ws.Range("H99").Value = 1428.5714
$ws.Range("I99").Value = 1300
$ws.Range("J99").Value = 1750
$ws.Range("K99").Value = 1300
$ws.Range("L99").Value = 1750
$ws.Range("M99").Value = 198
$ws.Range("N99").Value = -4746
$ws.Range("H105").Value = 2470.3
$ws.Range("I105").Value = 2213.125
$ws.Range("J105").Value = 3499
$ws.Range("K105").Value = 2213.125
$ws.Range("L105").Value = 3499
$ws.Range("M105").Value = -466.125
$ws.Range("N105").Value = -6993
$ws.Range("H107").Value = 72780.875
$ws.Range("I107").Value = 104442.2
$ws.Range("K107").Value = 104442.2
$ws.Range("M107").Value = -102522.2
$ws.Range("H134").Value = 2561.7896
$ws.Range("I134").Value = 2071.1333
$ws.Range("K134").Value = 6213.3999
$ws.Range("M134").Value = -3678.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1833.35
$ws.Range("I16").Value = 1526.7778
$ws.Range("J16").Value = 2084.182
$ws.Range("K16").Value = 1526.7778
$ws.Range("L16").Value = 2084.182
$ws.Range("M16").Value = -1239.7778
$ws.Range("N16").Value = -2658.182
$ws.Range("H22").Value = 1059.3
$ws.Range("I22").Value = 632.3333
$ws.Range("K22").Value = 632.3333
$ws.Range("M22").Value = -282.3333
$ws.Range("H31").Value = 3336235.8
$ws.Range("I31").Value = 6251115
$ws.Range("J31").Value = 4945.143
$ws.Range("K31").Value = 6251115
$ws.Range("L31").Value = 4945.143
$ws.Range("M31").Value = -6250820
$ws.Range("N31").Value = -5535.143
$ws.Range("H34").Value = 3336235.8
$ws.Range("I34").Value = 6251115
$ws.Range("J34").Value = 4945.143
$ws.Range("K34").Value = 6251115
$ws.Range("L34").Value = 4945.143
$ws.Range("M34").Value = -6250913
$ws.Range("N34").Value = -5349.143
$ws.Range("H51").Value = 36966
$ws.Range("J51").Value = 36966
$ws.Range("L51").Value = 36966
$ws.Range("N51").Value = -38438
$ws.Range("H58").Value = 1729.2333
$ws.Range("I58").Value = 1641.6818
$ws.Range("J58").Value = 1970
$ws.Range("K58").Value = 1641.6818
$ws.Range("L58").Value = 1970
$ws.Range("M58").Value = -1438.6818
$ws.Range("N58").Value = -2376
$ws.Range("H61").Value = 36966
$ws.Range("J61").Value = 36966
$ws.Range("L61").Value = 36966
$ws.Range("N61").Value = -37662
$ws.Range("H86").Value = 48946.81
$ws.Range("I86").Value = 71171.27
$ws.Range("J86").Value = 24499.9
$ws.Range("K86").Value = 71171.27
$ws.Range("L86").Value = 24499.9
$ws.Range("M86").Value = -70048.27
$ws.Range("N86").Value = -26745.9
$ws.Range("H89").Value = 48946.81
$ws.Range("I89").Value = 71171.27
$ws.Range("J89").Value = 24499.9
$ws.Range("K89").Value = 355856.35
$ws.Range("L89").Value = 122499.5
$ws.Range("M89").Value = -350240.35
$ws.Range("N89").Value = -133731.5
$ws.Range("H99").Value = 14411.333
$ws.Range("I99").Value = 10399.2
$ws.Range("J99").Value = 17277.143
$ws.Range("K99").Value = 10399.2
$ws.Range("L99").Value = 17277.143
$ws.Range("M99").Value = -8901.200000000001
$ws.Range("N99").Value = -20273.143
$ws.Range("H105").Value = 1817
$ws.Range("I105").Value = 1021.25
$ws.Range("K105").Value = 1021.25
$ws.Range("M105").Value = 725.75
$ws.Range("H107").Value = 898.86206
$ws.Range("I107").Value = 792.2727
$ws.Range("K107").Value = 792.2727
$ws.Range("M107").Value = 1127.7273
$ws.Range("H113").Value = 1833.35
$ws.Range("I113").Value = 1526.7778
$ws.Range("J113").Value = 2084.182
$ws.Range("K113").Value = 1526.7778
$ws.Range("L113").Value = 2084.182
$ws.Range("M113").Value = 643.2221999999999
$ws.Range("N113").Value = -6424.182
$ws.Range("H126").Value = 14411.333
$ws.Range("I126").Value = 10399.2
$ws.Range("J126").Value = 17277.143
$ws.Range("K126").Value = 31197.6
$ws.Range("L126").Value = 51831.429
$ws.Range("M126").Value = -28727.6
$ws.Range("N126").Value = -56771.429
$ws.Range("H132").Value = 78905
$ws.Range("I132").Value = 112207.22
$ws.Range("J132").Value = 3975
$ws.Range("K132").Value = 336621.66
$ws.Range("L132").Value = 11925
$ws.Range("M132").Value = -334091.66
$ws.Range("N132").Value = -16985
$ws.Range("H134").Value = 2041.8
$ws.Range("I134").Value = 2041.8
$ws.Range("K134").Value = 6125.4
$ws.Range("M134").Value = -3590.4
$ws.Range("H136").Value = 1729.2333
$ws.Range("I136").Value = 1641.6818
$ws.Range("J136").Value = 1970
$ws.Range("K136").Value = 4925.0454
$ws.Range("L136").Value = 5910
$ws.Range("M136").Value = -2375.0454
$ws.Range("N136").Value = -11010
$ws.Range("H141").Value = 117845.8
$ws.Range("J141").Value = 117845.8
$ws.Range("L141").Value = 117845.8
$ws.Range("N141").Value = -128205.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2049.5
$ws.Range("I5").Value = 2049.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 6148.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -6036.5
$ws.Range("N5").ClearContents()
$ws.Range("H12").Value = 148.04347
$ws.Range("J12").Value = 105.2
$ws.Range("L12").Value = 315.6
$ws.Range("N12").Value = -661.6
$ws.Range("H55").Value = 3565.6667
$ws.Range("I55").Value = 850
$ws.Range("K55").Value = 2550
$ws.Range("M55").Value = -2373
$ws.Range("H60").Value = 4366.6333
$ws.Range("I60").Value = 833.3333
$ws.Range("J60").Value = 4759.222
$ws.Range("K60").Value = 2499.9999
$ws.Range("L60").Value = 14277.666
$ws.Range("M60").Value = -2248.9999
$ws.Range("N60").Value = -14779.666
$ws.Range("H104").Value = 500
$ws.Range("I104").Value = 500
$ws.Range("K104").Value = 1500
$ws.Range("M104").Value = 1121
$ws.Range("H107").Value = 2103.6
$ws.Range("I107").Value = 4035.5
$ws.Range("J107").Value = 815.6667
$ws.Range("K107").Value = 12106.5
$ws.Range("L107").Value = 2447.0001
$ws.Range("M107").Value = -10186.5
$ws.Range("N107").Value = -6287.0001
$ws.Range("H110").Value = 11651.1
$ws.Range("I110").Value = 3302.2
$ws.Range("K110").Value = 9906.599999999999
$ws.Range("M110").Value = -5816.599999999999
$ws.Range("H122").Value = 1002.44446
$ws.Range("J122").Value = 1144.6154
$ws.Range("L122").Value = 10301.5386
$ws.Range("N122").Value = -15201.5386
$ws.Range("H130").Value = 7500
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H131").Value = 138659.22
$ws.Range("I131").Value = 251199.7
$ws.Range("J131").Value = 2002.9286
$ws.Range("K131").Value = 753599.1000000001
$ws.Range("L131").Value = 6008.7858
$ws.Range("M131").Value = -748559.1000000001
$ws.Range("N131").Value = -16088.7858
$ws.Range("H135").Value = 2049.5
$ws.Range("I135").Value = 2049.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 18445.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -15910.5
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 3083.4614
$ws.Range("I137").Value = 2400
$ws.Range("J137").Value = 3510.625
$ws.Range("K137").Value = 7200
$ws.Range("L137").Value = 10531.875
$ws.Range("M137").Value = -2100
$ws.Range("N137").Value = -20731.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12193.4375
$ws.Range("I70").Value = 13866.444
$ws.Range("J70").Value = 10042.429
$ws.Range("K70").Value = 13866.444
$ws.Range("L70").Value = 10042.429
$ws.Range("M70").Value = -13596.444
$ws.Range("N70").Value = -10582.429
$ws.Range("H73").Value = 12193.4375
$ws.Range("I73").Value = 13866.444
$ws.Range("J73").Value = 10042.429
$ws.Range("K73").Value = 13866.444
$ws.Range("L73").Value = 10042.429
$ws.Range("M73").Value = -12930.444
$ws.Range("N73").Value = -11914.429
$ws.Range("H80").Value = 10092.546
$ws.Range("I80").Value = 3839.8
$ws.Range("K80").Value = 3839.8
$ws.Range("M80").Value = -2841.8
$ws.Range("H83").Value = 10092.546
$ws.Range("I83").Value = 3839.8
$ws.Range("K83").Value = 19199
$ws.Range("M83").Value = -14207
$ws.Range("H99").Value = 8682.666999999999
$ws.Range("I99").Value = 8682.666999999999
$ws.Range("K99").Value = 8682.666999999999
$ws.Range("M99").Value = -6436.666999999999
$ws.Range("H100").Value = 23249
$ws.Range("J100").Value = 23249
$ws.Range("L100").Value = 23249
$ws.Range("N100").Value = -25413
$ws.Range("H102").Value = 5569.375
$ws.Range("I102").Value = 6960.273
$ws.Range("J102").Value = 2509.4
$ws.Range("K102").Value = 6960.273
$ws.Range("L102").Value = 2509.4
$ws.Range("M102").Value = -5338.273
$ws.Range("N102").Value = -5753.4
$ws.Range("H107").Value = 268.84616
$ws.Range("J107").Value = 371.33334
$ws.Range("L107").Value = 371.33334
$ws.Range("N107").Value = -4211.33334
$ws.Range("H122").Value = 2547.6897
$ws.Range("I122").Value = 2473.9583
$ws.Range("J122").Value = 2901.6
$ws.Range("K122").Value = 7421.874899999999
$ws.Range("L122").Value = 8704.799999999999
$ws.Range("M122").Value = -4971.874899999999
$ws.Range("N122").Value = -13604.8
$ws.Range("H126").Value = 2772.3076
$ws.Range("I126").Value = 1577.4286
$ws.Range("J126").Value = 4166.3335
$ws.Range("K126").Value = 4732.2858
$ws.Range("L126").Value = 12499.0005
$ws.Range("M126").Value = -2262.2858
$ws.Range("N126").Value = -17439.0005
$ws.Range("H132").Value = 2773.4412
$ws.Range("I132").Value = 2705.9697
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8117.909100000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5587.909100000001
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 905.56525
$ws.Range("I16").Value = 920.17645
$ws.Range("J16").Value = 864.1667
$ws.Range("K16").Value = 920.17645
$ws.Range("L16").Value = 864.1667
$ws.Range("M16").Value = -750.17645
$ws.Range("N16").Value = -1204.1667
$ws.Range("H17").Value = 56333
$ws.Range("J17").Value = 56333
$ws.Range("L17").Value = 56333
$ws.Range("N17").Value = -56673
$ws.Range("H40").Value = 2605.762
$ws.Range("I40").Value = 2151.4443
$ws.Range("J40").Value = 5331.6665
$ws.Range("K40").Value = 2151.4443
$ws.Range("L40").Value = 5331.6665
$ws.Range("M40").Value = -2015.4443
$ws.Range("N40").Value = -5603.6665
$ws.Range("H46").Value = 4531.1
$ws.Range("I46").Value = 712
$ws.Range("K46").Value = 712
$ws.Range("M46").Value = -524
$ws.Range("H61").Value = 1374.375
$ws.Range("I61").Value = 999.2857
$ws.Range("K61").Value = 999.2857
$ws.Range("M61").Value = -797.2857
$ws.Range("H82").Value = 2257.5715
$ws.Range("I82").Value = 2951
$ws.Range("J82").Value = 1980.2
$ws.Range("K82").Value = 2951
$ws.Range("L82").Value = 1980.2
$ws.Range("M82").Value = -2590
$ws.Range("N82").Value = -2702.2
$ws.Range("H85").Value = 2257.5715
$ws.Range("I85").Value = 2951
$ws.Range("J85").Value = 1980.2
$ws.Range("K85").Value = 2951
$ws.Range("L85").Value = 1980.2
$ws.Range("M85").Value = -1703
$ws.Range("N85").Value = -4476.2
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180
$ws.Range("H113").Value = 1374.375
$ws.Range("I113").Value = 999.2857
$ws.Range("K113").Value = 999.2857
$ws.Range("M113").Value = 1170.7143
$ws.Range("H122").Value = 4410.933
$ws.Range("I122").Value = 3805.4119
$ws.Range("J122").Value = 5202.769
$ws.Range("K122").Value = 11416.2357
$ws.Range("L122").Value = 15608.307
$ws.Range("M122").Value = -8966.235700000001
$ws.Range("N122").Value = -20508.307
$ws.Range("H132").Value = 3521.4285
$ws.Range("I132").Value = 3521.4285
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10564.2855
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8034.2855
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 4094.6
$ws.Range("I136").Value = 4105.1113
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 12315.3339
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -9765.333899999998
$ws.Range("N136").Value = -17100
$ws.Range("H138").Value = 99999.5
$ws.Range("J138").Value = 99999.5
$ws.Range("L138").Value = 99999.5
$ws.Range("N138").Value = -110279.5
$ws.Range("H141").Value = 110269
$ws.Range("J141").Value = 110269
$ws.Range("L141").Value = 110269
$ws.Range("N141").Value = -120629

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1439899.4
$ws.Range("J5").Value = 1439899.4
$ws.Range("L5").Value = 1439899.4
$ws.Range("N5").Value = -1440123.4
$ws.Range("H74").Value = 17468.2
$ws.Range("J74").Value = 19974
$ws.Range("L74").Value = 19974
$ws.Range("N74").Value = -21846
$ws.Range("H77").Value = 17468.2
$ws.Range("J77").Value = 19974
$ws.Range("L77").Value = 59922
$ws.Range("N77").Value = -69282
$ws.Range("H81").Value = 34916.168
$ws.Range("I81").Value = 40899.4
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 81798.8
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -80737.8
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 34916.168
$ws.Range("I84").Value = 40899.4
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 408994
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -403690
$ws.Range("N84").Value = -60608
$ws.Range("H100").Value = 1106.875
$ws.Range("I100").Value = 464.75
$ws.Range("K100").Value = 929.5
$ws.Range("M100").Value = -388.5
$ws.Range("H107").Value = 801.7
$ws.Range("I107").Value = 712.7619
$ws.Range("K107").Value = 2138.2857
$ws.Range("M107").Value = -218.2856999999999
$ws.Range("H113").Value = 970
$ws.Range("I113").Value = 800.1539
$ws.Range("J113").Value = 1154
$ws.Range("K113").Value = 2400.4617
$ws.Range("L113").Value = 3462
$ws.Range("M113").Value = -230.4616999999998
$ws.Range("N113").Value = -7802
$ws.Range("H122").Value = 36216.633
$ws.Range("I122").Value = 40782.816
$ws.Range("J122").Value = 6079.8
$ws.Range("K122").Value = 122348.448
$ws.Range("L122").Value = 18239.4
$ws.Range("M122").Value = -119898.448
$ws.Range("N122").Value = -23139.4
$ws.Range("H125").Value = 40749.75
$ws.Range("J125").Value = 40749.75
$ws.Range("L125").Value = 40749.75
$ws.Range("N125").Value = -50589.75
$ws.Range("H126").Value = 4646.5
$ws.Range("I126").Value = 3333
$ws.Range("J126").Value = 5004.727
$ws.Range("K126").Value = 9999
$ws.Range("L126").Value = 15014.181
$ws.Range("M126").Value = -7529
$ws.Range("N126").Value = -19954.181
$ws.Range("H132").Value = 35114.047
$ws.Range("I132").Value = 45206.312
$ws.Range("K132").Value = 135618.936
$ws.Range("M132").Value = -133088.936
$ws.Range("H136").Value = 22464.938
$ws.Range("I136").Value = 31144.137
$ws.Range("J136").Value = 3370.7
$ws.Range("K136").Value = 93432.41099999999
$ws.Range("L136").Value = 10112.1
$ws.Range("M136").Value = -90882.41099999999
$ws.Range("N136").Value = -15212.1
